$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 228
$ws.Range("I33").Value = 244.18182
$ws.Range("J33").Value = 50
$ws.Range("K33").Value = 244.18182
$ws.Range("L33").Value = 50
$ws.Range("M33").Value = -15.18181999999999
$ws.Range("N33").Value = -508
$ws.Range("H74").Value = 4811402.5
$ws.Range("I74").Value = 3299.5833
$ws.Range("J74").Value = 8932634
$ws.Range("K74").Value = 3299.5833
$ws.Range("L74").Value = 8932634
$ws.Range("M74").Value = -2363.5833
$ws.Range("N74").Value = -8934506
$ws.Range("H77").Value = 4811402.5
$ws.Range("I77").Value = 3299.5833
$ws.Range("J77").Value = 8932634
$ws.Range("K77").Value = 16497.9165
$ws.Range("L77").Value = 44663170
$ws.Range("M77").Value = -11817.9165
$ws.Range("N77").Value = -44672530
$ws.Range("H96").Value = 16667511
$ws.Range("J96").Value = 875.875
$ws.Range("L96").Value = 2627.625
$ws.Range("N96").Value = -5373.625
$ws.Range("H125").Value = 471.33334
$ws.Range("I125").Value = 380.16666
$ws.Range("K125").Value = 3421.49994
$ws.Range("M125").Value = -961.4999399999997
$ws.Range("H129").Value = 156983.56
$ws.Range("I129").Value = 280.16666
$ws.Range("K129").Value = 840.4999799999999
$ws.Range("M129").Value = 4159.50002
$ws.Range("H134").Value = 47970
$ws.Range("J134").Value = 47970
$ws.Range("L134").Value = 47970
$ws.Range("N134").Value = -58110

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1870.8649
$ws.Range("I61").Value = 1507.2258
$ws.Range("J61").Value = 3749.6667
$ws.Range("K61").Value = 1507.2258
$ws.Range("L61").Value = 3749.6667
$ws.Range("M61").Value = -1295.2258
$ws.Range("N61").Value = -4173.6667
$ws.Range("H97").Value = 830.3684
$ws.Range("I97").Value = 845.13336
$ws.Range("J97").Value = 775
$ws.Range("K97").Value = 845.13336
$ws.Range("L97").Value = 775
$ws.Range("M97").Value = -349.13336
$ws.Range("N97").Value = -1767
$ws.Range("H102").Value = 1446.1428
$ws.Range("I102").Value = 1353.8334
$ws.Range("K102").Value = 1353.8334
$ws.Range("M102").Value = 268.1666
$ws.Range("H124").Value = 11606.75
$ws.Range("J124").Value = 11606.75
$ws.Range("L124").Value = 11606.75
$ws.Range("N124").Value = -21426.75
$ws.Range("H132").Value = 11996.529
$ws.Range("I132").Value = 2046.0256
$ws.Range("K132").Value = 6138.0768
$ws.Range("M132").Value = -3608.0768
$ws.Range("H136").Value = 1870.8649
$ws.Range("I136").Value = 1507.2258
$ws.Range("J136").Value = 3749.6667
$ws.Range("K136").Value = 4521.6774
$ws.Range("L136").Value = 11249.0001
$ws.Range("M136").Value = -1971.6774
$ws.Range("N136").Value = -16349.0001

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1683.2354
$ws.Range("I20").Value = 1658.1111
$ws.Range("J20").Value = 1711.5
$ws.Range("K20").Value = 1658.1111
$ws.Range("L20").Value = 1711.5
$ws.Range("M20").Value = -1411.1111
$ws.Range("N20").Value = -2205.5
$ws.Range("H81").Value = 13008.637
$ws.Range("J81").Value = 13008.637
$ws.Range("L81").Value = 13008.637
$ws.Range("N81").Value = -15130.637
$ws.Range("H84").Value = 13008.637
$ws.Range("J84").Value = 13008.637
$ws.Range("L84").Value = 39025.911
$ws.Range("N84").Value = -49633.911
$ws.Range("H105").Value = 1088885
$ws.Range("I105").Value = 1640.45
$ws.Range("J105").Value = 1925226.9
$ws.Range("K105").Value = 1640.45
$ws.Range("L105").Value = 1925226.9
$ws.Range("M105").Value = 106.55
$ws.Range("N105").Value = -1928720.9
$ws.Range("H134").Value = 4500.1
$ws.Range("I134").Value = 4310.448
$ws.Range("K134").Value = 12931.344
$ws.Range("M134").Value = -10396.344

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 19578.822
$ws.Range("J58").Value = 34932.6
$ws.Range("L58").Value = 34932.6
$ws.Range("N58").Value = -35338.6
$ws.Range("H94").Value = 3025.5
$ws.Range("J94").Value = 3342.1667
$ws.Range("L94").Value = 3342.1667
$ws.Range("N94").Value = -4244.1667
$ws.Range("H99").Value = 3777.96
$ws.Range("I99").Value = 2713.1052
$ws.Range("K99").Value = 2713.1052
$ws.Range("M99").Value = -1215.1052
$ws.Range("H122").Value = 1055
$ws.Range("I122").Value = 907.9091
$ws.Range("K122").Value = 2723.7273
$ws.Range("M122").Value = -273.7273
$ws.Range("H126").Value = 3777.96
$ws.Range("I126").Value = 2713.1052
$ws.Range("K126").Value = 8139.3156
$ws.Range("M126").Value = -5669.3156
$ws.Range("H132").Value = 3237.8696
$ws.Range("I132").Value = 2310.875
$ws.Range("K132").Value = 6932.625
$ws.Range("M132").Value = -4402.625
$ws.Range("H136").Value = 19578.822
$ws.Range("J136").Value = 34932.6
$ws.Range("L136").Value = 104797.8
$ws.Range("N136").Value = -109897.8

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 275
$ws.Range("I7").Value = 50
$ws.Range("J7").Value = 500
$ws.Range("K7").Value = 150
$ws.Range("L7").Value = 1500
$ws.Range("M7").Value = -38
$ws.Range("N7").Value = -1724
$ws.Range("H80").Value = 19333.166
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 19333.166
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 57999.49800000001
$ws.Range("M80").ClearContents()
$ws.Range("N80").Value = -59871.49800000001
$ws.Range("H83").Value = 19333.166
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 19333.166
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 173998.494
$ws.Range("M83").ClearContents()
$ws.Range("N83").Value = -183358.494
$ws.Range("H92").Value = 470
$ws.Range("I92").Value = 100
$ws.Range("J92").Value = 562.5
$ws.Range("K92").Value = 300
$ws.Range("L92").Value = 1687.5
$ws.Range("M92").Value = 948
$ws.Range("N92").Value = -4183.5
$ws.Range("H131").Value = 744.45
$ws.Range("I131").Value = 291.2
$ws.Range("J131").Value = 768.30524
$ws.Range("K131").Value = 873.5999999999999
$ws.Range("L131").Value = 2304.91572
$ws.Range("M131").Value = 4166.4
$ws.Range("N131").Value = -12384.91572

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5694750
$ws.Range("I70").Value = 5250
$ws.Range("K70").Value = 5250
$ws.Range("M70").Value = -4980
$ws.Range("H73").Value = 5694750
$ws.Range("I73").Value = 5250
$ws.Range("K73").Value = 5250
$ws.Range("M73").Value = -4314
$ws.Range("H113").Value = 5504.591
$ws.Range("I113").Value = 8040.154
$ws.Range("J113").Value = 1842.1111
$ws.Range("K113").Value = 8040.154
$ws.Range("L113").Value = 1842.1111
$ws.Range("M113").Value = -5870.154
$ws.Range("N113").Value = -6182.1111

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2162.25
$ws.Range("I46").Value = 2212
$ws.Range("J46").Value = 2062.75
$ws.Range("K46").Value = 2212
$ws.Range("L46").Value = 2062.75
$ws.Range("M46").Value = -2024
$ws.Range("N46").Value = -2438.75
$ws.Range("H55").Value = 1146.8889
$ws.Range("I55").Value = 1404.2858
$ws.Range("J55").Value = 246
$ws.Range("K55").Value = 1404.2858
$ws.Range("L55").Value = 246
$ws.Range("M55").Value = -1231.2858
$ws.Range("N55").Value = -592
$ws.Range("H61").Value = 3739.0588
$ws.Range("I61").Value = 2018.8572
$ws.Range("K61").Value = 2018.8572
$ws.Range("M61").Value = -1816.8572
$ws.Range("H113").Value = 3739.0588
$ws.Range("I113").Value = 2018.8572
$ws.Range("K113").Value = 2018.8572
$ws.Range("M113").Value = 151.1428000000001
$ws.Range("H136").Value = 1467.3529
$ws.Range("I136").Value = 1363.9062
$ws.Range("J136").Value = 3122.5
$ws.Range("K136").Value = 4091.7186
$ws.Range("L136").Value = 9367.5
$ws.Range("M136").Value = -1541.7186
$ws.Range("N136").Value = -14467.5

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 39699.75
$ws.Range("J46").Value = 39699.75
$ws.Range("L46").Value = 39699.75
$ws.Range("N46").Value = -40161.75
$ws.Range("H100").Value = 259.3846
$ws.Range("I100").Value = 272.2
$ws.Range("J100").Value = 216.66667
$ws.Range("K100").Value = 544.4
$ws.Range("L100").Value = 433.33334
$ws.Range("M100").Value = -3.399999999999977
$ws.Range("N100").Value = -1515.33334
$ws.Range("H132").Value = 1119.8182
$ws.Range("I132").Value = 941.5
$ws.Range("J132").Value = 1431.875
$ws.Range("K132").Value = 2824.5
$ws.Range("L132").Value = 4295.625
$ws.Range("M132").Value = -294.5
$ws.Range("N132").Value = -9355.625
$ws.Range("H134").Value = 39699.75
$ws.Range("J134").Value = 39699.75
$ws.Range("L134").Value = 119099.25
$ws.Range("N134").Value = -124169.25
$ws.Range("H136").Value = 29496302
$ws.Range("I136").Value = 39703620
$ws.Range("J136").Value = 8499.444
$ws.Range("K136").Value = 119110860
$ws.Range("L136").Value = 25498.332
$ws.Range("M136").Value = -119108310
$ws.Range("N136").Value = -30598.332
